# Adds the new "Antenne GPS" order line to the "RS Component" sheet
# (row 7), wires up its RS Online hyperlink, and leaves that sheet/cell
# as the active selection - matching the author's commit
# "ajout antenne GPS aux commandes".

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("RS Component")

# --- new order line (row 7) ------------------------------------------------
$ws.Range("A7").Value = "Antenne GPS "
$ws.Range("B7").Value = "plus petit fil (50cm)"
$ws.Range("C7").Value = "140-8050 "
$ws.Range("D7").Value = "Siretta"
$ws.Range("E7").Value = "MIKE3A/0.5M/SMAM/RA/S/17"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 7.8
$ws.Range("H7").Value = "https://fr.rs-online.com/web/p/antennes-gps/1408050/"

# Turn H7 into a real hyperlink, same as the other rows in the table.
$ws.Hyperlinks.Add($ws.Range("H7"), "https://fr.rs-online.com/web/p/antennes-gps/1408050/")

# Hyperlinks.Add re-styles the cell (blue/underline); restore the plain
# "link-text" look used by the rest of column H (copy format only from H6).
$ws.Range("H6").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- active sheet / selection ----------------------------------------------
# The edit was finished on the "RS Component" sheet with F8 selected
# (previously "Mouser" was the active tab with D12 selected on RS Component).
[void]$ws.Activate()
$ws.Range("F8").Select() | Out-Null
